# Add a new "Save" column (H) to the s_vals sheet, matching the format
# of the existing header cells, with a numeric value of 0 in the data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (G1, which carries
# the bold/bordered/centered header style) onto the new header cell H1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Set the new header text and data value.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
